# Remove the bold "color red <comment>" placeholder runs that were
# inserted as review/feedback notes in front of several FirstParagraph
# paragraphs. Each target paragraph begins with one or more bold runs
# ("color", "red", then the actual comment text); everything else in
# the paragraph (if any, e.g. a trailing space + hyperlink) must be
# left untouched.

$d = $word.ActiveDocument

# Distinctive trailing text of each bold "color"/"red" comment run-group,
# used to locate the paragraph regardless of exact run segmentation.
$markers = @(
    "I looked this up and it seemed to be related to database management systems.",
    "before version detection and finally displaying to the user. Make a flowchart for this.",
    "Not really sure about this.",
    "This seems to be about the volume of data stored in a database."
)

$deletions = New-Object System.Collections.ArrayList

foreach ($m in $markers) {
    $found = $d.Content
    $ok = $found.Find.Execute($m, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "Marker not found (already removed?):" $m
        continue
    }

    # Expand to the whole enclosing paragraph so we know where the
    # bold comment block starts, even though Find only matched its tail.
    $para = $found.Duplicate
    [void]$para.Expand(4)   # wdParagraph

    $pStart = $para.Start
    $pEnd = $para.End - 1   # exclude the trailing paragraph mark

    # Walk forward from the paragraph start while characters are bold;
    # that prefix is exactly the "color"/"red"/comment run-group to drop.
    $pos = $pStart
    while ($pos -lt $pEnd) {
        $ch = $d.Range($pos, $pos + 1)
        if ($ch.Font.Bold -ne -1) {
            break
        }
        $pos = $pos + 1
    }

    if ($pos -gt $pStart) {
        [void]$deletions.Add(@{Start = $pStart; End = $pos})
    }
}

# Delete from the end of the document backwards so earlier ranges keep
# their character offsets valid.
$sorted = $deletions | Sort-Object -Property Start -Descending
foreach ($del in $sorted) {
    $delRange = $d.Range($del.Start, $del.End)
    Write-Host "Removing:" $delRange.Text
    $delRange.Delete()
}
